$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 contains sequential date serial numbers (A2:AD2) that need to be
# shifted forward by 13 days (e.g. A2: 45581 -> 45594, AD2: 45610 -> 45623).
for ($col = 1; $col -le 30; $col++) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.Value = $cell.Value2 + 13
}
